$excel = New-Object -ComObject Excel.Application
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '96.500.11'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.82%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.650.17'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.63%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '2.59'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +35.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.994'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.59%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '226.92'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.68%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '644.14'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.56%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.424'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.46%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.13'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +5.88%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.997'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.27%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '3.644.35'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.71%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '48.61'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +9.48%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.211'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.84%  '
$ws.Range('E14').Value = '  -7.99%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.63'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.67%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.310.58'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.05%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '95.658.21'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.47%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '21.46'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +14.64%  '
$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '8.87'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.61%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.20'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +8.08%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '3.635.84'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.535'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +6.07%  '
$ws.Range('B23').Value = 'Hedera'
$ws.Range('C23').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.264'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +35.30%  '
$ws.Range('B24').Value = 'BitcoinCash'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '517.14'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.82%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.28'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.41%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '122.62'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +19.64%  '
$ws.Range('E27').Value = '  -6.56%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.84'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.91%  '
$ws.Range('B29').Value = 'Aptos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '12.94'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.25%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '13.25'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.69%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.02'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.05%  '
$ws.Range('B32').Value = 'Dai'
$ws.Range('C32').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.00'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.05%  '
$ws.Range('B33').Value = 'Cronos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.184'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.55%  '
$ws.Range('B34').Value = 'PolygonEcosystemToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.626'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.17%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '33.12'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.62%  '
$ws.Range('B36').Value = 'Binance-PegBSC-USD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.996'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.90%  '
$ws.Range('B37').Value = 'Fetch.AI'
$ws.Range('C37').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.78'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.01%  '
$ws.Range('B38').Value = 'Bittensor'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '607.38'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -7.50%  '
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.48'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.16%  '
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.10'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.26%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '42.80'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +6.51%  '
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.496'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.62%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0502'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +9.08%  '
$ws.Range('B45').Value = 'Kaspa'
$ws.Range('C45').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.160'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.90%  '
$ws.Range('B46').Value = 'ImmutableX'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.97'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.960'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.34%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.31'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.01%  '
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.83'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.67%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '226.22'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +9.54%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '23.55'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.28%  '
